$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the existing values from column A into column B, and put new
# localized row labels into column A (row 3 is brand new).
$ws.Range("B1").Value = 0.9
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 1

$ws.Range("A1").Value = "Коэффициент скорости обучения"
$ws.Range("A2").Value = "Количество входов нейронной сети"
$ws.Range("A3").Value = "Размерность выходного слоя"

# Move the active selection to match the saved view state.
$ws.Range("N8").Select()
